# Add the new "Save" column (H) with header + data, matching the style of
# the existing header row (e.g. B1:G1) and plain numeric data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - copy style from the existing header G1 ("sum") so it
# gets the same bold/bordered/centered style (s="1"), then set its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column, rows 2-13.
$saveValues = @(0, 0, 0, 1, 1, 1, 0, 1, 1, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
